# feat: add 2022-Q4 data
#
# 1. Insert a brand-new "2022-Q4" sheet right after "总计" (and before the
#    existing "2021-Q1" sheet), populated with the new quarter's fund data.
# 2. Update the "总计" (totals) sheet: the former row-2 ("2021-Q1") becomes
#    the new "2022-Q4" summary row, the former row-3 ("2020-Q4") becomes
#    "2021-Q1" (unchanged numbers), and a brand-new row-4 is appended for
#    "2020-Q4" (repeating the same holding numbers).

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q4" worksheet right after "总计" (sheet 1)
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

# Pull in the header/index-column formatting (style "2") used by every
# other quarter sheet, by copying formats from the "总计" sheet, which
# already carries that style on its header row + index column.
$total.Range("B1").Copy() | Out-Null
$q4.Range("B1:H1").PasteSpecial($xlPasteFormats) | Out-Null
$total.Range("A2").Copy() | Out-Null
$q4.Range("A2:A6").PasteSpecial($xlPasteFormats) | Out-Null

# Header row
$q4.Cells.Item(1,2).Value = "基金代码"
$q4.Cells.Item(1,3).Value = "基金名称"
$q4.Cells.Item(1,4).Value = "基金规模"
$q4.Cells.Item(1,5).Value = "股票总仓位"
$q4.Cells.Item(1,6).Value = "仓位占比"
$q4.Cells.Item(1,7).Value = "持有市值(亿元)"
$q4.Cells.Item(1,8).Value = "仓位排名"

# Columns B, D, E, F, G hold numeric-looking data (fund codes with
# leading zeros, percentages, scale, etc.) that must stay TEXT - format
# as text first so Excel doesn't silently coerce them into numbers.
# (Column C is the fund name, never numeric-looking, so it's left alone.)
$q4.Range("B2:B6").NumberFormat = "@"
$q4.Range("D2:G6").NumberFormat = "@"

$q4Rows = @(
    @{ idx = 0; code = "012466"; name = "嘉实策略精选混合A";       size = "9.59"; pos = "93.89"; pct = "2.90"; mv = "0.2781"; rank = 10 },
    @{ idx = 1; code = "513690"; name = "博时恒生港股通高股息率ETF"; size = "5.20"; pos = "98.04"; pct = "2.66"; mv = "0.1383"; rank = 7 },
    @{ idx = 2; code = "014307"; name = "嘉实多元动力混合A";       size = "1.56"; pos = "93.22"; pct = "3.20"; mv = "0.0499"; rank = 8 },
    @{ idx = 3; code = "012467"; name = "嘉实策略精选混合C";       size = "1.01"; pos = "93.89"; pct = "2.90"; mv = "0.0293"; rank = 10 },
    @{ idx = 4; code = "014308"; name = "嘉实多元动力混合C";       size = "0.33"; pos = "93.22"; pct = "3.20"; mv = "0.0106"; rank = 8 }
)

$r = 2
foreach ($row in $q4Rows) {
    $q4.Cells.Item($r,1).Value = $row.idx
    $q4.Cells.Item($r,2).Value = $row.code
    $q4.Cells.Item($r,3).Value = $row.name
    $q4.Cells.Item($r,4).Value = $row.size
    $q4.Cells.Item($r,5).Value = $row.pos
    $q4.Cells.Item($r,6).Value = $row.pct
    $q4.Cells.Item($r,7).Value = $row.mv
    $q4.Cells.Item($r,8).Value = $row.rank
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) Update the "总计" sheet
# ---------------------------------------------------------------------

# Row 2 : was "2021-Q1" data -> becomes the new "2022-Q4" summary
$total.Cells.Item(2,2).Value = "2022-Q4"
$total.Cells.Item(2,3).Value = 5
$total.Cells.Item(2,4).Value = 0.51

# Row 3 : was "2020-Q4" data -> becomes "2021-Q1" (same counts as before)
$total.Cells.Item(3,2).Value = "2021-Q1"
$total.Cells.Item(3,3).Value = 2
$total.Cells.Item(3,4).Value = 0.01

# Row 4 (new) : "2020-Q4", repeating the same holding numbers
$total.Cells.Item(4,2).Value = "2020-Q4"
$total.Cells.Item(4,3).Value = 2
$total.Cells.Item(4,4).Value = 0.01
$total.Cells.Item(4,1).Value = 2
$total.Range("A3").Copy() | Out-Null
$total.Range("A4").PasteSpecial($xlPasteFormats) | Out-Null

# Adding a worksheet makes it the active tab; restore "总计" as the
# active/selected sheet (the workbook's bookViews were untouched by
# this change, so the original sheet stays the one in focus).
$total.Activate()
